# Atualização de bases das ligas, do dia: 28-05-2024 às 20:56
#
# Two pairs of adjacent match rows had their data entered in the wrong
# order (the two fixtures on the same matchday got swapped). Fix this by
# swapping the contents of columns B:AD (everything except the running
# index in column A) between each pair of rows:
#   rows 110 <-> 111
#   rows 192 <-> 193

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Costa Rica Primera Division")

function Swap-RowData($row1, $row2) {
    $range1 = $ws.Range("B$row1`:AD$row1")
    $range2 = $ws.Range("B$row2`:AD$row2")

    $values1 = $range1.Value2
    $values2 = $range2.Value2

    $range1.Value = $values2
    $range2.Value = $values1
}

Swap-RowData 110 111
Swap-RowData 192 193
